# Design Specifications and Budget updated
# Modified Design Specifications to be submitted September 23

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Budget sheet: add the "so far" running total to row 9 ---
# B9 gets a bold, currency-formatted value of 131 ...
$ws.Range("B9").Value = 131
$ws.Range("B9").Font.Bold = $true
$ws.Range("B9").NumberFormat = """$""#,##0.00"

# ... and C9 gets the note "so far"
$ws.Range("C9").Value = "so far"

# The active selection moved from A9 to B9
[void]$ws.Range("B9").Select()

# --- Title textbox: clarify the project name ---
$shp = $ws.Shapes.Item("TextBox 1")
$shp.TextFrame.Characters().Text = "EE4951W - Power Monitoring Device - Initial Budget`n`nOnly add total system cost, not individual items."
